$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A4").Value = "hgjkhk"
$ws.Range("B4").Value = "gffgh"
$ws.Range("A5").Value = "65h889"
$ws.Range("B5").Value = 778

$ws.Range("B5").Select()
